# VEG_to_FUEL_TYPE.xlsx — "fixes and input layers available"
#
# Fix the "shurbs" -> "shrubs" typo wherever it appears in the fuel-type
# column (B), and move the selection to where the author was last working
# (D29) as captured by the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B28").Value = "shrubs"
$ws.Range("B29").Value = "shrubs"
$ws.Range("B30").Value = "shrubs"
$ws.Range("B33").Value = "shrubs"
$ws.Range("B34").Value = "shrubs"

$ws.Range("D29").Select() | Out-Null
